$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.03125
$ws.Range("E2").Value = 0.134
$ws.Range("G2").Value = 0.3926582278481013
$ws.Range("H2").Value = 0.3926582278481013
$ws.Range("I2").Value = 0.4139240506329114
$ws.Range("J2").Value = 0.2986642119648333
$ws.Range("K2").Value = 11.64
$ws.Range("L2").Value = 0.2946835443037975
$ws.Range("M2").Value = 4.7028
$ws.Range("N2").Value = 0.02926446795270691
$ws.Range("O2").Value = 0.404020618556701
$ws.Range("P2").Value = 4.7028
$ws.Range("Q2").Value = 0.02926446795270691
$ws.Range("R2").Value = 0.404020618556701
$ws.Range("U2").Value = 36.99
$ws.Range("V2").Value = 0.2301804604853765
$ws.Range("W2").Value = 0.1052374893977947
$ws.Range("X2").Value = 0.06353781670736119
$ws.Range("Y2").Value = 0.04169967269043355
$ws.Range("Z2").Value = 0.9139287366959741
$ws.Range("AA2").Value = 0.5586622342402616
$ws.Range("AB2").Value = 0.06353781670736119
$ws.Range("AC2").Value = 0.4947803635985127
$ws.Range("AD2").Value = 0.952
$ws.Range("AF2").Value = 0.952
$ws.Range("AG2").Value = -36.038
$ws.Range("AH2").Value = 0.005889194071214708
$ws.Range("AI2").Value = 0.008084788368775053
$ws.Range("AJ2").Value = -0.2890856876995396
$ws.Range("AK2").Value = -0.4462247096406726
$ws.Range("AL2").Value = 0.113
$ws.Range("AM2").Value = 0.113
$ws.Range("AN2").Value = 0.05477560414269275
$ws.Range("AO2").Value = 144.6902654867257
$ws.Range("AP2").Value = -2.073532796317607
$ws.Range("AQ2").Value = 144.6902654867257

# Row 3
$ws.Range("B3").Value = "Peoples Insurance Company Limited (DSE:PEOPLESINS)"
$ws.Range("D3").Value = 0.0382
$ws.Range("E3").ClearContents()
$ws.Range("G3").Value = 0.2264957264957265
$ws.Range("H3").Value = 0.2264957264957265
$ws.Range("I3").Value = 0.2663817663817664
$ws.Range("J3").Value = 0.194017094017094
$ws.Range("K3").Value = 1.36
$ws.Range("L3").Value = 0.1937321937321937
$ws.Range("M3").Value = 0.4158
$ws.Range("N3").Value = 0.01506521739130435
$ws.Range("O3").Value = 0.305735294117647
$ws.Range("P3").Value = 0.4158
$ws.Range("Q3").Value = 0.01506521739130435
$ws.Range("R3").Value = 0.305735294117647
$ws.Range("U3").Value = 15.6
$ws.Range("V3").Value = 0.5652173913043478
$ws.Range("W3").Value = 0.09444444444444444
$ws.Range("X3").Value = 0.06353781670736119
$ws.Range("Y3").Value = 0.03090662773708325
$ws.Range("Z3").Value = 7.799999999999996
$ws.Range("AA3").Value = 1.513333333333333
$ws.Range("AB3").Value = 0.06353781670736119
$ws.Range("AC3").Value = 1.449795516625971
$ws.Range("AG3").Value = -15.6
$ws.Range("AJ3").Value = -1.3
$ws.Range("AK3").Value = -156.0000000000005
$ws.Range("AP3").Value = -7.393364928909953

# Row 4
$ws.Range("B4").Value = "Bangladesh General Insurance Company Limited (DSE:BGIC)"
$ws.Range("D4").Value = 0.0243
$ws.Range("E4").Value = 0.134
$ws.Range("G4").Value = 0.3376623376623377
$ws.Range("H4").Value = 0.3376623376623377
$ws.Range("I4").Value = 0.3327922077922078
$ws.Range("J4").Value = 0.2638071563852813
$ws.Range("K4").Value = 1.52
$ws.Range("L4").Value = 0.2467532467532468
$ws.Range("M4").Value = 0.701
$ws.Range("N4").Value = 0.02577205882352941
$ws.Range("O4").Value = 0.4611842105263158
$ws.Range("P4").Value = 0.701
$ws.Range("Q4").Value = 0.02577205882352941
$ws.Range("R4").Value = 0.4611842105263158
$ws.Range("U4").Value = 11
$ws.Range("V4").Value = 0.4044117647058824
$ws.Range("W4").Value = 0.116030534351145
$ws.Range("X4").Value = 0.06496155939592943
$ws.Range("Y4").Value = 0.05106897495521562
$ws.Range("Z4").Value = 3.242105263157894
$ws.Range("AA4").Value = 0.8552905701754384
$ws.Range("AB4").Value = 0.06422592457613664
$ws.Range("AC4").Value = 0.7910646455993018
$ws.Range("AD4").Value = 0.952
$ws.Range("AF4").Value = 0.952
$ws.Range("AG4").Value = -10.048
$ws.Range("AH4").Value = 0.03381642512077294
$ws.Range("AI4").Value = 0.072939013178057
$ws.Range("AJ4").Value = -0.585820895522388
$ws.Range("AK4").Value = -4.896686159844055
$ws.Range("AL4").Value = 0.113
$ws.Range("AM4").Value = 0.113
$ws.Range("AN4").Value = 0.4387096774193548
$ws.Range("AO4").Value = 18.14159292035398
$ws.Range("AP4").Value = -4.630414746543779
$ws.Range("AQ4").Value = 18.14159292035398

# Row 5
$ws.Range("D5").Value = 0.0491
$ws.Range("E5").Value = 0.153
$ws.Range("G5").Value = 0.4794520547945206
$ws.Range("H5").Value = 0.4794520547945206
$ws.Range("I5").Value = 0.4931506849315069
$ws.Range("J5").Value = 0.3529680365296805
$ws.Range("K5").Value = 7.71
$ws.Range("L5").Value = 0.3520547945205479
$ws.Range("M5").Value = 3.03
$ws.Range("N5").Value = 0.03935064935064935
$ws.Range("O5").Value = 0.3929961089494163
$ws.Range("P5").Value = 3.03
$ws.Range("Q5").Value = 0.03935064935064935
$ws.Range("R5").Value = 0.3929961089494163
$ws.Range("U5").Value = 3.08
$ws.Range("V5").Value = 0.04
$ws.Range("W5").Value = 0.1189814814814815
$ws.Range("X5").Value = 0.06353781670736119
$ws.Range("Y5").Value = 0.05544366477412029
$ws.Range("Z5").Value = 0.7423728813559322
$ws.Range("AA5").Value = 0.2620338983050848
$ws.Range("AB5").Value = 0.06353781670736119
$ws.Range("AC5").Value = 0.1984960815977236
$ws.Range("AG5").Value = -3.08
$ws.Range("AJ5").Value = -0.04166666666666666
$ws.Range("AK5").Value = -0.04534746760895171
$ws.Range("AP5").Value = -0.2725663716814159

# Row 6
$ws.Range("B6").Value = "United Insurance Company Limited (DSE:UNITEDINS)"
$ws.Range("D6").Value = 0.0072
$ws.Range("E6").Value = -0.0488
$ws.Range("G6").Value = 0.3031674208144797
$ws.Range("H6").Value = 0.3031674208144797
$ws.Range("I6").Value = 0.3687782805429864
$ws.Range("J6").Value = 0.2394782414390258
$ws.Range("K6").Value = 1.05
$ws.Range("L6").Value = 0.2375565610859729
$ws.Range("M6").Value = 0.556
$ws.Range("N6").Value = 0.01923875432525952
$ws.Range("O6").Value = 0.5295238095238095
$ws.Range("P6").Value = 0.556
$ws.Range("Q6").Value = 0.01923875432525952
$ws.Range("R6").Value = 0.5295238095238095
$ws.Range("U6").Value = 7.31
$ws.Range("V6").Value = 0.2529411764705882
$ws.Range("W6").Value = 0.062874251497006
$ws.Range("X6").Value = 0.06353781670736119
$ws.Range("Y6").Value = -0.0006635652103551964
$ws.Range("Z6").Value = 0.4047619047619048
$ws.Range("AA6").Value = 0.09693166915389141
$ws.Range("AB6").Value = 0.06353781670736119
$ws.Range("AC6").Value = 0.03339385244653022
$ws.Range("AD6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = -7.31
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = -0.3385826771653543
$ws.Range("AK6").Value = -0.6838166510757716
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").Value = 0
$ws.Range("AO6").ClearContents()
$ws.Range("AP6").Value = -4.061111111111111
$ws.Range("AQ6").ClearContents()
